# [IMP] New test data
# Adds payment_term_id values for the rows that were still missing them
# (rows 4, 5, 7 and 8 of the sale_order test sheet) and nudges a couple of
# cosmetic view settings (tab ratio / selected cell) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data in column H (payment_term_id) ---------------------------
# H4 and H8 reuse the same value/style as the existing H1/H2 cells
# (z0bug.payment_1), H5 and H7 get new payment refs and reuse H3's style.

$ws.Range("H2").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H4").Value = "z0bug.payment_1"

$ws.Range("H3").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H5").Value = "z0bug.payment_5"

$ws.Range("H3").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H7").Value = "z0bug.payment_4"

$ws.Range("H2").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H8").Value = "z0bug.payment_1"

$excel.CutCopyMode = 0

# --- View tweaks --------------------------------------------------------
$wb.Windows.Item(1).TabRatio = 0.5
$ws.Range("H7").Select() | Out-Null
